# Applies the "angled_lens" model re-train update:
#  - appends 32 new prediction rows (154-185) from the refreshed model run
#  - keeps the `_xlchart.v1.*` defined-name bookkeeping in sync with what
#    Excel re-derives for the box/whisker charts after the sheet refresh
#  - leaves the sheet scrolled/selected near the newly added rows

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("angled_lens")

# ---------------------------------------------------------------------
# 1. Append the 32 new rows of data (columns A:M) below the existing 153
# ---------------------------------------------------------------------
$newRows = @(
    @(1, 1, -0.13738962267853999, 0, 1, 0.51587703602726598, 0, 1, 0, 0, 0, 0, "right"),
    @(0, 0, [double]"-3.9210615080581997E-2", 1, 0, -0.70358898427225403, 1, 0, 1, 1, 1, 0, "left"),
    @(1, 1, -0.102619998531325, 0, 1, 0.58974518722783797, 0, 1, 0, 0, 0, 0, "right"),
    @(1, 1, 0.206725430578036, 1, 0, -0.53174863780797899, 0, 0, 0, 1, 1, 1, "left"),
    @(1, 0, 0.48798652063946701, 0, 0, 0.50624694473689502, 0, 1, 1, 1, 0, 0, "right"),
    @(0, 1, -0.37806983269879302, 1, 0, -0.44788265580733, 1, 0, 0, 1, 1, 1, "left"),
    @(1, 0, 0.194314234221142, 0, 1, 0.58985349585669, 0, 1, 1, 0, 0, 0, "right"),
    @(1, 0, 0.24350109352591601, 1, 1, -0.134318306610581, 0, 0, 1, 0, 1, 1, "left"),
    @(1, 0, 0.37594547062094302, 0, 1, 0.568360709737739, 0, 1, 1, 1, 0, 0, "right"),
    @(1, 1, [double]"4.5140250812914999E-2", 1, 1, [double]"-7.3240283100208997E-2", 0, 0, 0, 0, 1, 1, "left"),
    @(1, 1, -0.16375735357659399, 0, 1, 0.57558489305098903, 0, 1, 0, 0, 0, 0, "right"),
    @(1, 1, [double]"-7.5502637399739003E-2", 1, 1, -0.152070072121571, 0, 0, 0, 0, 1, 1, "left"),
    @(1, 0, 0.21027304138696001, 0, 1, 0.52947422451341197, 0, 1, 1, 1, 0, 0, "right"),
    @(1, 0, 0.21049130578338901, 1, 1, -0.107316276816038, 0, 0, 1, 0, 1, 1, "left"),
    @(1, 0, 0.156062314922704, 0, 1, 0.57589516671274399, 0, 1, 1, 1, 0, 0, "right"),
    @(1, 0, 0.22530120363701001, 1, 1, [double]"-5.8773997350353997E-2", 0, 0, 1, 0, 1, 1, "left"),
    @(1, 0, 0.20975386895383999, 0, 1, 0.59143480040240204, 0, 1, 1, 1, 0, 0, "right"),
    @(0, 0, [double]"-5.3589796034666001E-2", 1, 0, -0.44116138640098101, 1, 0, 1, 1, 1, 1, "left"),
    @(1, 0, 0.118097000032258, 0, 1, 0.58355679818873396, 0, 1, 1, 0, 0, 0, "right"),
    @(1, 1, [double]"2.2126371731535002E-2", 1, 0, -0.43297523314370501, 0, 0, 0, 1, 1, 1, "left"),
    @(1, 0, [double]"-7.8127783411249992E-3", 0, 1, 0.57549723120140195, 0, 1, 1, 0, 0, 0, "right"),
    @(1, 0, 0.24714067545398299, 1, 1, [double]"-8.5041133398212995E-2", 0, 0, 1, 0, 1, 1, "left"),
    @(1, 1, [double]"-3.6960445064643001E-2", 0, 1, 0.51996780701656298, 0, 1, 0, 1, 0, 0, "right"),
    @(0, 0, [double]"-3.1647155362772003E-2", 1, 1, [double]"-8.8591282741013994E-2", 1, 0, 1, 0, 1, 1, "left"),
    @(1, 1, -0.113467651190855, 0, 1, 0.59880642242594895, 0, 1, 0, 0, 0, 0, "right"),
    @(0, 0, [double]"6.3032824387702996E-2", 1, 0, -0.45504249910063499, 1, 0, 1, 1, 1, 1, "left"),
    @(1, 1, -0.13805112357704599, 0, 1, 0.58491461213679996, 0, 1, 0, 1, 0, 0, "right"),
    @(0, 0, [double]"2.5131818367733001E-2", 1, 0, -0.430603666451496, 1, 0, 1, 1, 1, 0, "left"),
    @(1, 1, [double]"-4.6062546737821002E-2", 0, 1, 0.55551737989622096, 0, 1, 0, 0, 0, 0, "right"),
    @(1, 0, 0.380606181479971, 1, 0, -0.588463870852306, 0, 0, 1, 1, 1, 0, "left"),
    @(1, 0, 0.21261656138578, 0, 1, 0.55384132410582898, 0, 1, 1, 0, 0, 0, "right"),
    @(1, 1, 0.14534180474835701, 1, 1, [double]"-5.5894741287923003E-2", 0, 0, 1, 1, 1, 1, "left")
)

$rowCount = $newRows.Count
$colCount = 13
$block = New-Object 'object[,]' $rowCount,$colCount
for ($i = 0; $i -lt $rowCount; $i++) {
    for ($j = 0; $j -lt $colCount; $j++) {
        $block[$i,$j] = $newRows[$i][$j]
    }
}

$startRow = 154
$endRow = $startRow + $rowCount - 1
$targetRange = $ws.Range("A" + $startRow + ":M" + $endRow)
$targetRange.Value = $block

# ---------------------------------------------------------------------
# 2. Re-sync the hidden `_xlchart.v1.*` defined names used by the
#    box-and-whisker charts on wall_mounted_data (Excel reassigns these
#    IDs whenever the backing charts are refreshed)
# ---------------------------------------------------------------------
$names = $wb.Names
$definedNameMap = @{
    "_xlchart.v1.1"  = "wall_mounted_data!`$B`$1";
    "_xlchart.v1.2"  = "wall_mounted_data!`$B`$2:`$B`$320";
    "_xlchart.v1.4"  = "wall_mounted_data!`$C`$1";
    "_xlchart.v1.5"  = "wall_mounted_data!`$C`$2:`$C`$320";
    "_xlchart.v1.7"  = "wall_mounted_data!`$D`$1";
    "_xlchart.v1.8"  = "wall_mounted_data!`$D`$2:`$D`$320";
    "_xlchart.v1.10" = "wall_mounted_data!`$F`$1";
    "_xlchart.v1.11" = "wall_mounted_data!`$F`$2:`$F`$320";
    "_xlchart.v1.13" = "wall_mounted_data!`$E`$1";
    "_xlchart.v1.14" = "wall_mounted_data!`$E`$2:`$E`$320";
}
foreach ($key in $definedNameMap.Keys) {
    $names.Item($key).RefersTo = "=" + $definedNameMap[$key]
}

# ---------------------------------------------------------------------
# 3. Leave the sheet active, scrolled to the freshly added rows, with the
#    selection resting near the new data (no more pinned top-left cell)
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("P15").Select()

Write-Output ("angled_lens used range: " + $ws.UsedRange.Address())
